$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 17 (shifts existing rows 17..61 down to 18..62)
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new data record
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44519
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 100112022
$ws.Cells.Item(17, 7).Value = "Arveja Verde"
$ws.Cells.Item(17, 8).Value = "Perfection"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 65
$ws.Cells.Item(17, 11).Value = 15000
$ws.Cells.Item(17, 12).Value = 16000
$ws.Cells.Item(17, 13).Value = 15538
$ws.Cells.Item(17, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Región del Maule"
$ws.Cells.Item(17, 16).Value = 622
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"
